# Update cryptocurrency price/volume data (and reorder a few rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.240.92"
$ws.Range("E2").Value = "'  +0.70%  "

$ws.Range("D3").Value = "'3.926.79"
$ws.Range("E3").Value = "'  +1.88%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  -0.03%  "

$ws.Range("D5").Value = "'483.36"
$ws.Range("E5").Value = "'  +5.79%  "

$ws.Range("D6").Value = "'147.96"
$ws.Range("E6").Value = "'  +0.75%  "

$ws.Range("D7").Value = "'0.622"
$ws.Range("E7").Value = "'  -0.19%  "

$ws.Range("D8").Value = "'0.998"
$ws.Range("E8").Value = "'  -0.10%  "

$ws.Range("D9").Value = "'0.727"
$ws.Range("E9").Value = "'  -2.43%  "

$ws.Range("D10").Value = "'0.167"
$ws.Range("E10").Value = "'  +7.49%  "

$ws.Range("D11").Value = "'0.0000354"
$ws.Range("E11").Value = "'  +11.01%  "

$ws.Range("D12").Value = "'42.62"
$ws.Range("E12").Value = "'  -2.52%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'10.43"
$ws.Range("E13").Value = "'  +0.81%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "'4.530.03"
$ws.Range("E14").Value = "'  +1.49%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "'4.009.75"
$ws.Range("E15").Value = "'  +4.43%  "

$ws.Range("B16").Value = "Uniswap"
$ws.Range("C16").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D16").Value = "'14.72"
$ws.Range("E16").Value = "'  -1.06%  "

$ws.Range("E17").Value = "'  -0.26%  "

$ws.Range("D18").Value = "'19.88"
$ws.Range("E18").Value = "'  -1.24%  "

$ws.Range("D19").Value = "'1.14"
$ws.Range("E19").Value = "'  -2.72%  "

$ws.Range("D20").Value = "'68.151.84"
$ws.Range("E20").Value = "'  +0.43%  "

$ws.Range("D21").Value = "'434.47"
$ws.Range("E21").Value = "'  +1.31%  "

$ws.Range("D22").Value = "'3.42"
$ws.Range("E22").Value = "'  +5.87%  "

$ws.Range("D23").Value = "'14.48"
$ws.Range("E23").Value = "'  -2.09%  "

$ws.Range("D24").Value = "'87.27"
$ws.Range("E24").Value = "'  +0.61%  "

$ws.Range("D25").Value = "'10.84"
$ws.Range("E25").Value = "'  +5.73%  "

$ws.Range("D26").Value = "'3.57"
$ws.Range("E26").Value = "'  +1.01%  "

$ws.Range("D27").Value = "'38.39"
$ws.Range("E27").Value = "'  +2.65%  "

$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").Value = "'10.67"
$ws.Range("E28").Value = "'  +8.74%  "

$ws.Range("B29").Value = "LEO"
$ws.Range("C29").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D29").Value = "'5.89"
$ws.Range("E29").Value = "'  +7.83%  "

$ws.Range("D30").Value = "'720.28"
$ws.Range("E30").Value = "'  -2.77%  "

$ws.Range("D31").Value = "'13.30"
$ws.Range("E31").Value = "'  -3.33%  "

$ws.Range("D32").Value = "'0.129"
$ws.Range("E32").Value = "'  -3.54%  "

$ws.Range("D33").Value = "'2.83"
$ws.Range("E33").Value = "'  +3.08%  "

$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").Value = "'41.98"
$ws.Range("E34").Value = "'  -3.00%  "

$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D35").Value = "'0.0₃0889"
$ws.Range("E35").Value = "'  +28.84%  "

$ws.Range("D36").Value = "'59.07"
$ws.Range("E36").Value = "'  +2.75%  "

$ws.Range("E37").Value = "'  -6.32%  "

$ws.Range("D38").Value = "'5.51"
$ws.Range("E38").Value = "'  -1.24%  "

$ws.Range("E39").Value = "'  -0.18%  "

$ws.Range("D40").Value = "'2.89"
$ws.Range("E40").Value = "'  +8.68%  "

$ws.Range("D41").Value = "'0.0470"
$ws.Range("E41").Value = "'  -1.27%  "

$ws.Range("E42").Value = "'  +11.26%  "

$ws.Range("D43").Value = "'2.98"
$ws.Range("E43").Value = "'  +1.01%  "

$ws.Range("D44").Value = "'0.348"
$ws.Range("E44").Value = "'  -1.52%  "

$ws.Range("D45").Value = "'0.141"
$ws.Range("E45").Value = "'  +0.74%  "

$ws.Range("D46").Value = "'0.999"
$ws.Range("E46").Value = "'  -0.25%  "

$ws.Range("E47").Value = "'  -0.19%  "

$ws.Range("D48").Value = "'2.18"
$ws.Range("E48").Value = "'  +1.90%  "

$ws.Range("D49").Value = "'3.25"
$ws.Range("E49").Value = "'  -3.03%  "

$ws.Range("D50").Value = "'145.59"
$ws.Range("E50").Value = "'  +1.29%  "

$ws.Range("D51").Value = "'2.85"
$ws.Range("E51").Value = "'  -0.99%  "
